$d = $word.ActiveDocument

# --- Change 1: "Email excel file to Maryam " paragraph ---
# Split the single run "Email excel file to Maryam " into two runs:
#   "Email excel file to Maryam" (no trailing space) + " - DONE"
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Email excel file to Maryam") {
        $hit = $p.Range.Duplicate
        $hit.Find.Execute("Email excel file to Maryam ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        $target = $d.Range($hit.Start, $hit.End)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Email excel file to Maryam</w:t></w:r><w:r><w:t xml:space="preserve"> - DONE</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $target.InsertXML($xml)
        break
    }
}

# --- Change 2: "Remove mrn,surgery_date" paragraph ---
# Append a new run " - DONE" right after the existing "mrn,surgery_date" run
# (which sits right after a </w:proofErr>, so a plain InsertAfter on a
# collapsed range at paragraph-end naturally lands in its own new run).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^Remove mrn") {
        $r = $p.Range
        $endPoint = $d.Range($r.End - 1, $r.End - 1)
        $endPoint.InsertAfter(" - DONE")
        break
    }
}

# --- Change 3: add two more empty paragraphs before the trailing empty one ---
$last = $d.Paragraphs.Last
$r = $last.Range
$target = $d.Range($r.Start, $r.End)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
